# "arreglo de funciones para lectura de archivo excel"
#
# The validation sheet marks rows found in "Satelite" with the shared
# string "Existe en Satelite" in column C. The lookup/read logic was
# fixed so that rows which do NOT match now correctly get flagged with a
# new label "No Existe en Satelite" instead of being left with the
# (wrong) "Existe en Satelite" value. The header row (C1) and the first
# data row (C2) are untouched - only the actual validation rows
# (C3:C7) are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Validación ventas redenciones R")

$ws.Range("C3:C7").Value = "No Existe en Satelite"

# Reflect where the user was working when the fix was made: column C
# (the validation column) selected top to bottom.
$ws.Range("C1:C1048576").Select()
